$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89. This pushes the existing rows 89..195
# down to 90..196, preserving all of their data and formatting, and
# leaves a blank (but correctly formatted) row 89 for us to populate.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the missing "Especial" quality
# record for the 2021-08-04 Comercializadora del Agro de Limari entry.
$ws.Range("A89").Value = 2
$ws.Range("B89").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C89").Value = "Coquimbo"
$ws.Range("D89").Value = 44412
$ws.Range("E89").Value = 4
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100101
$ws.Range("H89").Value = "Berries"
$ws.Range("I89").Value = 100112025
$ws.Range("J89").Value = "Frutilla"
$ws.Range("K89").Value = "Sin especificar"
$ws.Range("L89").Value = "Especial"
$ws.Range("M89").Value = 240
$ws.Range("N89").Value = 25500
$ws.Range("O89").Value = 26000
$ws.Range("P89").Value = 25750
$ws.Range("Q89").Value = "$/bandeja 7 kilos"
$ws.Range("R89").Value = "Provincia de Melipilla"
$ws.Range("S89").Value = 3679
$ws.Range("T89").Value = 7
